# Run prepare & render with final data.
# Adds Russia as a new country column (inserted before Saudi Arabia), refreshes the
# aggregate ("All"/"Europe") computations and per-country survey values with the
# final data pull, and rewords the "petition" survey question / rewraps the
# "Governments..." question text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L (Saudi Arabia), shifting Saudi Arabia -> M
# and USA -> N, to make room for the new Russia column.
$ws.Columns.Item(12).Insert()

$ws.Cells.Item(1,1).Value = ""
$ws.Cells.Item(1,2).Value = "`$ bold('All')"
$ws.Cells.Item(1,3).Value = "`$ bold('Europe')"
$ws.Cells.Item(1,4).Value = "France"
$ws.Cells.Item(1,5).Value = "Germany"
$ws.Cells.Item(1,6).Value = "Italy"
$ws.Cells.Item(1,7).Value = "Poland"
$ws.Cells.Item(1,8).Value = "Spain"
$ws.Cells.Item(1,9).Value = "United Kingdom"
$ws.Cells.Item(1,10).Value = "Switzerland"
$ws.Cells.Item(1,11).Value = "Japan"
$ws.Cells.Item(1,12).Value = "Russia"
$ws.Cells.Item(1,13).Value = "Saudi Arabia"
$ws.Cells.Item(1,14).Value = "USA"
$ws.Cells.Item(2,1).Value = "Supports tax on world top 1% to finance global poverty reduction`n(Additional 15% tax on income over [`$120k/year in PPP])"
$ws.Cells.Item(2,2).Value = 0.556007222541788
$ws.Cells.Item(2,3).Value = 0.612606004275779
$ws.Cells.Item(2,4).Value = 0.616724928283738
$ws.Cells.Item(2,5).Value = 0.620808072627315
$ws.Cells.Item(2,6).Value = 0.75031578480093
$ws.Cells.Item(2,7).Value = 0.495076291993524
$ws.Cells.Item(2,8).Value = 0.6117343914093
$ws.Cells.Item(2,9).Value = 0.554433842033875
$ws.Cells.Item(2,10).Value = 0.529872342725065
$ws.Cells.Item(2,11).Value = 0.438282499020205
$ws.Cells.Item(2,12).Value = 0.595756191953926
$ws.Cells.Item(2,13).Value = 0.677827112481047
$ws.Cells.Item(2,14).Value = 0.508231644030169
$ws.Cells.Item(3,1).Value = "Supports tax on world top 3% to finance global poverty reduction`n(Additional 15% tax over [`$80k], 30% over [`$120k], 45% over [`$1M])"
$ws.Cells.Item(3,2).Value = 0.50312291438834
$ws.Cells.Item(3,3).Value = 0.55532380171867
$ws.Cells.Item(3,4).Value = 0.588306427984715
$ws.Cells.Item(3,5).Value = 0.528389866344171
$ws.Cells.Item(3,6).Value = 0.602852192238351
$ws.Cells.Item(3,7).Value = 0.547884560154944
$ws.Cells.Item(3,8).Value = 0.567216620094722
$ws.Cells.Item(3,9).Value = 0.542572051573976
$ws.Cells.Item(3,10).Value = 0.360111742646164
$ws.Cells.Item(3,11).Value = 0.351050882580874
$ws.Cells.Item(3,12).Value = 0.604399737467109
$ws.Cells.Item(3,13).Value = 0.671355992475847
$ws.Cells.Item(3,14).Value = 0.447300643788012
$ws.Cells.Item(4,1).Value = "Prefers sustainable future"
$ws.Cells.Item(4,2).Value = 0.680881448179833
$ws.Cells.Item(4,3).Value = 0.701539116816613
$ws.Cells.Item(4,4).Value = 0.719673932778947
$ws.Cells.Item(4,5).Value = 0.704645953003574
$ws.Cells.Item(4,6).Value = 0.760434225154626
$ws.Cells.Item(4,7).Value = 0.572739428737827
$ws.Cells.Item(4,8).Value = 0.73608837047601
$ws.Cells.Item(4,9).Value = 0.679200342355269
$ws.Cells.Item(4,10).Value = 0.665468662233116
$ws.Cells.Item(4,11).Value = 0.758076861129753
$ws.Cells.Item(4,12).Value = 0.688615273248795
$ws.Cells.Item(4,13).Value = 0.713280127381035
$ws.Cells.Item(4,14).Value = 0.616918649447641
$ws.Cells.Item(5,1).Value = "`"Governments should actively cooperate to have all countries`nconverge in terms of GDP per capita by the end of the century`""
$ws.Cells.Item(5,2).Value = 0.609601586795904
$ws.Cells.Item(5,3).Value = 0.667387097439935
$ws.Cells.Item(5,4).Value = 0.608837961112165
$ws.Cells.Item(5,5).Value = 0.64586557994442
$ws.Cells.Item(5,6).Value = 0.780212171052138
$ws.Cells.Item(5,7).Value = 0.766346455024569
$ws.Cells.Item(5,8).Value = 0.728439095410156
$ws.Cells.Item(5,9).Value = 0.555251737340332
$ws.Cells.Item(5,10).Value = 0.566109070039295
$ws.Cells.Item(5,11).Value = 0.571096670838126
$ws.Cells.Item(5,12).Value = 0.778963825426238
$ws.Cells.Item(5,13).Value = 0.756844881931732
$ws.Cells.Item(5,14).Value = 0.472900191628792
$ws.Cells.Item(6,1).Value = "Would support a global movement to tackle CC, tax millionaires,`n and fund LICs (either petition, demonstrate, strike, or donate)"
$ws.Cells.Item(6,2).Value = 0.675595447215337
$ws.Cells.Item(6,3).Value = 0.719216740354837
$ws.Cells.Item(6,4).Value = 0.699222514786681
$ws.Cells.Item(6,5).Value = 0.688082663981164
$ws.Cells.Item(6,6).Value = 0.819338712934373
$ws.Cells.Item(6,7).Value = 0.708420268414952
$ws.Cells.Item(6,8).Value = 0.74352822863702
$ws.Cells.Item(6,9).Value = 0.68138828161491
$ws.Cells.Item(6,10).Value = 0.639531813440066
$ws.Cells.Item(6,11).Value = 0.557841849059486
$ws.Cells.Item(6,13).Value = 0.727098526374066
$ws.Cells.Item(6,14).Value = 0.666651932459956
$ws.Cells.Item(7,1).Value = "More likely to vote for party if part of worldwide`ncoalition for climate action and global redistribution"
$ws.Cells.Item(7,2).Value = 0.364717906507653
$ws.Cells.Item(7,3).Value = 0.40721776360235
$ws.Cells.Item(7,4).Value = 0.428789858236926
$ws.Cells.Item(7,5).Value = 0.382760045850686
$ws.Cells.Item(7,6).Value = 0.49409402735291
$ws.Cells.Item(7,7).Value = 0.297020517023446
$ws.Cells.Item(7,8).Value = 0.434344229636586
$ws.Cells.Item(7,9).Value = 0.394744985134441
$ws.Cells.Item(7,10).Value = 0.314122252360097
$ws.Cells.Item(7,11).Value = 0.22156020948145
$ws.Cells.Item(7,14).Value = 0.372217577193357
$ws.Cells.Item(8,1).Value = "Supports reparations for colonization and slavery in`nthe form of funding education and technology transfers"
$ws.Cells.Item(8,2).Value = 0.347853243460036
$ws.Cells.Item(8,3).Value = 0.389632627264691
$ws.Cells.Item(8,4).Value = 0.330491389745475
$ws.Cells.Item(8,5).Value = 0.335223833737993
$ws.Cells.Item(8,6).Value = 0.55884551964724
$ws.Cells.Item(8,8).Value = 0.402288413192778
$ws.Cells.Item(8,9).Value = 0.353117864669717
$ws.Cells.Item(8,14).Value = 0.306445646731996
$ws.Cells.Item(9,1).Value = "`"My taxes should go towards solving global problems`""
$ws.Cells.Item(9,2).Value = 0.410626908494325
$ws.Cells.Item(9,3).Value = 0.439007996245569
$ws.Cells.Item(9,4).Value = 0.324278224348557
$ws.Cells.Item(9,5).Value = 0.441488389669955
$ws.Cells.Item(9,6).Value = 0.518375799884677
$ws.Cells.Item(9,7).Value = 0.40468182659417
$ws.Cells.Item(9,8).Value = 0.506211274976579
$ws.Cells.Item(9,9).Value = 0.451738941455142
$ws.Cells.Item(9,10).Value = 0.393105560133332
$ws.Cells.Item(9,11).Value = 0.319383802321488
$ws.Cells.Item(9,12).Value = 0.389233362357354
$ws.Cells.Item(9,13).Value = 0.662536527281334
$ws.Cells.Item(9,14).Value = 0.40055514051731
$ws.Cells.Item(10,1).Value = "`"My taxes ... global problems`" (Global Nation, 2024)"
$ws.Cells.Item(10,2).Value = 0.0
$ws.Cells.Item(10,3).Value = 0.0
$ws.Cells.Item(10,4).Value = 0.0
$ws.Cells.Item(10,5).Value = 0.0
$ws.Cells.Item(10,6).Value = 0.0
$ws.Cells.Item(10,7).Value = 0.0
$ws.Cells.Item(10,8).Value = 0.0
$ws.Cells.Item(10,9).Value = 0.0
$ws.Cells.Item(10,11).Value = 0.0
$ws.Cells.Item(10,13).Value = 0.0
$ws.Cells.Item(10,14).Value = 0.0
